$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 264; this shifts the existing rows 264-372 down to 265-373,
# matching the target layout (old row N -> new row N+1 for N in 264..372).
$ws.Rows("264:264").Insert()

# Populate the newly inserted row 264 with the new weekly record.
$ws.Range("A264").Value = 5
$ws.Range("B264").Value = "Macroferia Regional de Talca"
$ws.Range("C264").Value = "Maule"
$ws.Range("D264").Value = 45027
$ws.Range("D264").NumberFormat = $ws.Range("D265").NumberFormat
$ws.Range("E264").Value = 7
$ws.Range("F264").Value = "Fruta"
$ws.Range("G264").Value = 100108
$ws.Range("H264").Value = "Tropicales y subtropicales"
$ws.Range("I264").Value = 100108005
$ws.Range("J264").Value = "Piña"
$ws.Range("K264").Value = "Caramelo"
$ws.Range("L264").Value = "Segunda"
$ws.Range("M264").Value = 200
$ws.Range("N264").Value = 19000
$ws.Range("O264").Value = 19000
$ws.Range("P264").Value = 19000
$ws.Range("Q264").Value = '$/caja 14 unidades'
$ws.Range("R264").Value = "Ecuador"
$ws.Range("S264").Value = 1357
$ws.Range("T264").Value = 14
